# Fix "reasource_req" typo -> "resource_req", and add the empty-requirement
# entry for House (G2) so the "can afford to build" check has something to
# read for every row.

$wb = $excel.ActiveWorkbook

$editable  = $wb.Worksheets.Item("editable")
$formatted = $wb.Worksheets.Item("formatted")

# Fix the misspelled header.
$editable.Range("G1").Value = "resource_req"

# House has no resource requirement to build -> empty JSON object.
$editable.Range("G2").Value = "{}"

# Selections / active tab: "editable" becomes the active sheet/tab now
# (it was "formatted" before). Set the no-longer-active sheet's selection
# first, then finish on "editable" so it is left as the active tab.
$formatted.Activate()
$formatted.Range("A2:H4").Select()

$editable.Activate()
$editable.Range("E14").Select()
